# Apply updated "想去人数" (F) / "最低票价" (G) figures to both the
# "展览" and "全部类型" sheets, which hold duplicate event listings.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 141
    $ws.Range("G2").Value = 60

    $ws.Range("F3").Value = 1713

    $ws.Range("G4").Value = 55

    $ws.Range("F9").Value = 637
}
